$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain decimal-looking numbers (e.g. "6.52", "18.09").
# The sheet stores them as text, so force a Text number format on those specific
# cells before assigning, otherwise Excel would silently convert them to numbers.
$textCells = @("D5", "D6", "D10", "D16", "D19", "D20", "D21", "D24", "D26", "D27", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D43", "D44", "D45", "D46", "D48", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '59.252.84'
$ws.Range("E2").Value = '  +1.82%  '
$ws.Range("D3").Value = '2.588.13'
$ws.Range("E3").Value = '  -0.32%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '523.57'
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").Value = '139.42'
$ws.Range("E6").Value = '  -2.98%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("E8").Value = '  -0.74%  '
$ws.Range("D9").Value = '2.599.08'
$ws.Range("E9").Value = '  -0.70%  '
$ws.Range("D10").Value = '6.52'
$ws.Range("E10").Value = '  -1.68%  '
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("E12").Value = '  -2.08%  '
$ws.Range("E13").Value = '  +2.83%  '
$ws.Range("D14").Value = '3.046.19'
$ws.Range("E14").Value = '  -0.35%  '
$ws.Range("D15").Value = '59.013.27'
$ws.Range("E15").Value = '  +1.45%  '
$ws.Range("D16").Value = '20.51'
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("D18").Value = '2.586.00'
$ws.Range("E18").Value = '  -1.14%  '
$ws.Range("D19").Value = '341.81'
$ws.Range("E19").Value = '  +0.73%  '
$ws.Range("D20").Value = '4.31'
$ws.Range("E20").Value = '  -1.31%  '
$ws.Range("D21").Value = '10.08'
$ws.Range("E21").Value = '  -2.05%  '
$ws.Range("E22").Value = '  +1.02%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").Value = '66.60'
$ws.Range("E24").Value = '  +1.97%  '
$ws.Range("E25").Value = '  +0.77%  '
$ws.Range("D26").Value = '0.406'
$ws.Range("E26").Value = '  +0.75%  '
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("E28").Value = '  +0.51%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").Value = '0.0₃0724'
$ws.Range("E30").Value = '  -3.20%  '
$ws.Range("D31").Value = '5.90'
$ws.Range("E31").Value = '  -5.29%  '
$ws.Range("D32").Value = '1.59'
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").Value = '18.72'
$ws.Range("E33").Value = '  -0.40%  '
$ws.Range("D34").Value = '149.26'
$ws.Range("E34").Value = '  -0.31%  '
$ws.Range("D35").Value = '3.98'
$ws.Range("E35").Value = '  -1.34%  '
$ws.Range("E36").Value = '  -1.81%  '
$ws.Range("D37").Value = '36.78'
$ws.Range("E37").Value = '  +2.10%  '
$ws.Range("D38").Value = '1.48'
$ws.Range("E38").Value = '  +1.35%  '
$ws.Range("D39").Value = '0.827'
$ws.Range("E39").Value = '  -3.80%  '
$ws.Range("D40").Value = '0.816'
$ws.Range("E40").Value = '  -6.71%  '
$ws.Range("E41").Value = '  -0.68%  '
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("D43").Value = '272.09'
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("D44").Value = '0.598'
$ws.Range("E44").Value = '  -0.49%  '
$ws.Range("D45").Value = '10.77'
$ws.Range("E45").Value = '  +0.89%  '
$ws.Range("D46").Value = '0.0951'
$ws.Range("E46").Value = '  -0.77%  '
$ws.Range("E47").Value = '  -1.74%  '
$ws.Range("D48").Value = '18.41'
$ws.Range("E48").Value = '  -2.26%  '
$ws.Range("D49").Value = '1.970.62'
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '18.09'
$ws.Range("E51").Value = '  -4.56%  '
